# Commit: update file with jgit
# Change cell E8 on the "Rules" sheet from "Good Morning" to "GIT UPDATE",
# and leave the selection on E8 (as the diff shows E8 becoming the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Select()
$ws.Range("E8").Value = "GIT UPDATE"
